$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B2 holds the Neo4j "Case ID" list query (WebExcel / query tab).
# It gains a trailing "order By ... LIMIT 100" clause, replacing the tab
# before "demo.survival_time" with a plain newline.
$newQuery = @'
MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
WHERE  samp.tissue_type in ["Tumor"]
return DISTINCT ss.study_subject_id as `Case ID`,
   p.program_acronym as `Program Code`,
    p.program_id as Program_ID,
   s.study_acronym as `Arm`,
   ss.disease_subtype as `Diagnosis`,
   sf.grouped_recurrence_score AS `Recurrence Score`,
   d.tumor_size_group AS `tumor_size`,
   d.er_status AS `ER Status`,
   d.pr_status AS `PR Status`,
   demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
order By ss.study_subject_id ASC LIMIT 100 
'@

# Strip the single trailing newline the here-string literal adds, keeping
# the single trailing space that is part of the real content.
$newQuery = $newQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newQuery

# Reflect the updated selection (B2 is now the active cell, was C2).
[void]$ws.Range("B2").Select()
